$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Posicao") previously held textual position names (e.g. "Extremo
# Esquerdo"). The edit replaces those labels with a numeric position code
# (1-13), one per row, following the mapping derived from the workbook.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 4
$ws.Range("C8").Value = 5
$ws.Range("C9").Value = 3
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 3
$ws.Range("C14").Value = 3
$ws.Range("C15").Value = 3
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 2
$ws.Range("C20").Value = 7
$ws.Range("C21").Value = 6
$ws.Range("C22").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("C24").Value = 4
$ws.Range("C25").Value = 3
$ws.Range("C26").Value = 3
$ws.Range("C27").Value = 3
$ws.Range("C28").Value = 6
$ws.Range("C29").Value = 6
$ws.Range("C30").Value = 8
$ws.Range("C31").Value = 6
$ws.Range("C32").Value = 6
$ws.Range("C33").Value = 4
$ws.Range("C34").Value = 1
$ws.Range("C35").Value = 2
$ws.Range("C36").Value = 3
$ws.Range("C37").Value = 6
$ws.Range("C38").Value = 1
$ws.Range("C39").Value = 1
$ws.Range("C40").Value = 3
$ws.Range("C41").Value = 4
$ws.Range("C42").Value = 3
$ws.Range("C43").Value = 6
$ws.Range("C44").Value = 7
$ws.Range("C45").Value = 1
$ws.Range("C46").Value = 3
$ws.Range("C47").Value = 9
$ws.Range("C48").Value = 10
$ws.Range("C49").Value = 1
$ws.Range("C50").Value = 9
$ws.Range("C51").Value = 8
$ws.Range("C52").Value = 10
$ws.Range("C53").Value = 8
$ws.Range("C54").Value = 11
$ws.Range("C55").Value = 6
$ws.Range("C56").Value = 8
$ws.Range("C57").Value = 3
$ws.Range("C58").Value = 3
$ws.Range("C59").Value = 11
$ws.Range("C60").Value = 9
$ws.Range("C61").Value = 7
$ws.Range("C62").Value = 9
$ws.Range("C63").Value = 10
$ws.Range("C64").Value = 3
$ws.Range("C65").Value = 5
$ws.Range("C66").Value = 3
$ws.Range("C67").Value = 1
$ws.Range("C68").Value = 4
$ws.Range("C69").Value = 9
$ws.Range("C70").Value = 9
$ws.Range("C71").Value = 2
$ws.Range("C72").Value = 9
$ws.Range("C73").Value = 5
$ws.Range("C74").Value = 10
$ws.Range("C75").Value = 9
$ws.Range("C76").Value = 11
$ws.Range("C77").Value = 11
$ws.Range("C78").Value = 4
$ws.Range("C79").Value = 6
$ws.Range("C80").Value = 9
$ws.Range("C81").Value = 9
$ws.Range("C82").Value = 3
$ws.Range("C83").Value = 1
$ws.Range("C84").Value = 2
$ws.Range("C85").Value = 2
$ws.Range("C86").Value = 8
$ws.Range("C87").Value = 6
$ws.Range("C88").Value = 3
$ws.Range("C89").Value = 8
$ws.Range("C90").Value = 6
$ws.Range("C91").Value = 9
$ws.Range("C92").Value = 9
$ws.Range("C93").Value = 8
$ws.Range("C94").Value = 2
$ws.Range("C95").Value = 7
$ws.Range("C96").Value = 1
$ws.Range("C97").Value = 9
$ws.Range("C98").Value = 2
$ws.Range("C99").Value = 9
$ws.Range("C100").Value = 10
$ws.Range("C101").Value = 6
$ws.Range("C102").Value = 7
$ws.Range("C103").Value = 3
$ws.Range("C104").Value = 7
$ws.Range("C105").Value = 3
$ws.Range("C106").Value = 11
$ws.Range("C107").Value = 6
$ws.Range("C108").Value = 9
$ws.Range("C109").Value = 7
$ws.Range("C110").Value = 9
$ws.Range("C111").Value = 1
$ws.Range("C112").Value = 2
$ws.Range("C113").Value = 7
$ws.Range("C114").Value = 9
$ws.Range("C115").Value = 6
$ws.Range("C116").Value = 6
$ws.Range("C117").Value = 2
$ws.Range("C118").Value = 7
$ws.Range("C119").Value = 9
$ws.Range("C120").Value = 1
$ws.Range("C121").Value = 6
$ws.Range("C122").Value = 6
$ws.Range("C123").Value = 1
$ws.Range("C124").Value = 2
$ws.Range("C125").Value = 9
$ws.Range("C126").Value = 9
$ws.Range("C127").Value = 9
$ws.Range("C128").Value = 1
$ws.Range("C129").Value = 6
$ws.Range("C130").Value = 9
$ws.Range("C131").Value = 11
$ws.Range("C132").Value = 2
$ws.Range("C133").Value = 7
$ws.Range("C134").Value = 9
$ws.Range("C135").Value = 9
$ws.Range("C136").Value = 1
$ws.Range("C137").Value = 4
$ws.Range("C138").Value = 9
$ws.Range("C139").Value = 9
$ws.Range("C140").Value = 9
$ws.Range("C141").Value = 6
$ws.Range("C142").Value = 1
$ws.Range("C143").Value = 8
$ws.Range("C144").Value = 6
$ws.Range("C145").Value = 6
$ws.Range("C146").Value = 9
$ws.Range("C147").Value = 3
$ws.Range("C148").Value = 2
$ws.Range("C149").Value = 9
$ws.Range("C150").Value = 12
$ws.Range("C151").Value = 11
$ws.Range("C152").Value = 2
$ws.Range("C153").Value = 2
$ws.Range("C154").Value = 1
$ws.Range("C155").Value = 10
$ws.Range("C156").Value = 11
$ws.Range("C157").Value = 7
$ws.Range("C158").Value = 9
$ws.Range("C159").Value = 3
$ws.Range("C160").Value = 6
$ws.Range("C161").Value = 1
$ws.Range("C162").Value = 6
$ws.Range("C163").Value = 4
$ws.Range("C164").Value = 7
$ws.Range("C165").Value = 7
$ws.Range("C166").Value = 6
$ws.Range("C167").Value = 6
$ws.Range("C168").Value = 11
$ws.Range("C169").Value = 6
$ws.Range("C170").Value = 4
$ws.Range("C171").Value = 4
$ws.Range("C172").Value = 6
$ws.Range("C173").Value = 2
$ws.Range("C174").Value = 1
$ws.Range("C175").Value = 7
$ws.Range("C176").Value = 2
$ws.Range("C177").Value = 9
$ws.Range("C178").Value = 2
$ws.Range("C179").Value = 9
$ws.Range("C180").Value = 6
$ws.Range("C181").Value = 4
$ws.Range("C182").Value = 3
$ws.Range("C183").Value = 6
$ws.Range("C184").Value = 13
$ws.Range("C185").Value = 4
$ws.Range("C186").Value = 4
$ws.Range("C187").Value = 7
$ws.Range("C188").Value = 9
$ws.Range("C189").Value = 9
$ws.Range("C190").Value = 9
$ws.Range("C191").Value = 7
$ws.Range("C192").Value = 4
$ws.Range("C193").Value = 6
$ws.Range("C194").Value = 10
$ws.Range("C195").Value = 6
$ws.Range("C196").Value = 2
$ws.Range("C197").Value = 3
$ws.Range("C198").Value = 9
$ws.Range("C199").Value = 9
$ws.Range("C200").Value = 10
$ws.Range("C201").Value = 6

# Resize column C to fit its new (numeric) content and select C3, matching
# the saved view state of the edited workbook.
$ws.Columns.Item(3).ColumnWidth = 17.1666666666667
$ws.Range("C3").Select() | Out-Null
